$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("I2:I10").Value = "d"
$ws.Range("A1:K1").Font.Bold = $true
$ws.Range("I2:I10").Select() | Out-Null
